$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Authors (E2) now point to a fresh/updated author-vote string, the
# old Misc. Data (I2) tag is cleared, and a new "Other found locations" (J2)
# value records the publisher.
$ws.Range("E2").Value = "[Wei%Tang%NULL%0, Zhujun%Cao%NULL%1, Mingfeng%Han%NULL%0, Zhengyan%Wang%NULL%1, Junwen%Chen%NULL%1, Wenjin%Sun%NULL%1, Yaojie%Wu%NULL%1, Wei%Xiao%NULL%0, Shengyong%Liu%NULL%1, Erzhen%Chen%NULL%1, Wei%Chen%NULL%0, Xiongbiao%Wang%NULL%1, Jiuyong%Yang%NULL%1, Jun%Lin%NULL%0, Qingxia%Zhao%NULL%1, Youqin%Yan%NULL%0, Zhibin%Xie%NULL%1, Dan%Li%NULL%0, Yaofeng%Yang%NULL%1, Leshan%Liu%NULL%1, Jieming%Qu%NULL%0, Guang%Ning%NULL%1, Guochao%Shi%NULL%1, Qing%Xie%NULL%1]"
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = "BMJ Publishing Group Ltd."

# Row 3: this citation could no longer be resolved with full text, so its
# Title/Abstract/Authors/ID/ID Format/Date Accepted/Misc. Data all revert to
# the "unknown" placeholders.
$ws.Range("C3").Value = "Unknown Title"
$ws.Range("D3").Value = "Unknown Abstract"
$ws.Range("E3").Value = "[]"
$ws.Range("F3").Value = "not found"
$ws.Range("G3").Value = "N/A"
# Force this as literal text (not an auto-converted date serial number), then
# drop the number-format override so the cell keeps the sheet's plain style.
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "1970-01-01"
$ws.Range("H3").ClearFormats()
$ws.Range("I3").Value = ""

# Row 4: same fallback-to-unknown treatment as row 3 (Date Accepted was
# already 1970-01-01, so it is left untouched).
$ws.Range("C4").Value = "Unknown Title"
$ws.Range("D4").Value = "Unknown Abstract"
$ws.Range("E4").Value = "[]"
$ws.Range("F4").Value = "not found"
$ws.Range("G4").Value = "N/A"
$ws.Range("I4").Value = ""

# Row 5: Authors (E5) refreshed with updated vote flags, Misc. Data (I5)
# cleared, and a new "Other found locations" (J5) value records the publisher.
$ws.Range("E5").Value = "[Philippe%Gautret%NULL%0, Jean-Christophe%Lagier%NULL%1, Philippe%Parola%NULL%0, Van Thuan%Hoang%NULL%0, Line%Meddeb%NULL%1, Morgane%Mailhe%NULL%1, Barbara%Doudier%NULL%1, Johan%Courjon%NULL%1, Valérie%Giordanengo%NULL%1, Vera Esteves%Vieira%NULL%1, Hervé%Tissot Dupont%NULL%1, Stéphane%Honoré%NULL%1, Philippe%Colson%NULL%1, Eric%Chabrière%NULL%1, Bernard%La Scola%NULL%1, Jean-Marc%Rolain%NULL%1, Philippe%Brouqui%NULL%1, Didier%Raoult%Didier.raoult@gmail.com%0]"
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = "Published by Elsevier B.V."
